$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duplicate row 2's formatting down into the two new rows (3 and 4) ---
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F4").PasteSpecial(-4122)

# --- Clear the stale "Results" (PASS/FAIL) values before retyping text ---
$ws.Cells.Item(2,5).ClearContents()
$ws.Cells.Item(3,5).ClearContents()
$ws.Cells.Item(4,5).ClearContents()

# --- Row 3 (new): Transmittals_New_ConsultantAdvice (typed first) ---
$ws.Cells.Item(3,1).Value = "Transmittals_New_ConsultantAdvice"
$ws.Cells.Item(3,1).WrapText = $true

# --- Row 2: Transmittals_New -> Transmittals_New_Correspondence ---
$ws.Cells.Item(2,1).Value = "Transmittals_New_Correspondence"
$ws.Cells.Item(2,1).WrapText = $true
$ws.Cells.Item(2,2).Value = "Creates a new Transmittal of Type Correspondence"
$ws.Cells.Item(2,3).Value = "N"
$ws.Cells.Item(2,4).Value = "Y"
$ws.Cells.Item(2,6).Value = "Sprint1"

# --- finish row 3 ---
$ws.Cells.Item(3,2).Value = "Creates a new Transmittal of Type Consultant Advice"
$ws.Cells.Item(3,3).Value = "N"
$ws.Cells.Item(3,4).Value = "Y"
$ws.Cells.Item(3,6).Value = "Sprint1"

# --- Row 4 (new): Transmittals_New_ChangeNote ---
$ws.Cells.Item(4,1).Value = "Transmittals_New_ChangeNote"
$ws.Cells.Item(4,1).WrapText = $true
$ws.Cells.Item(4,2).Value = "Creates a new Transmittal of Type Change Note"
$ws.Cells.Item(4,3).Value = "N"
$ws.Cells.Item(4,4).Value = "Y"
$ws.Cells.Item(4,6).Value = "Sprint1"

# --- Data validations: extend list ranges to cover the new rows ---
$rngYN = $ws.Range("C2:D4")
$rngYN.Validation.Delete()
$rngYN.Validation.Add(3, 1, 1, '"Y,N"')

$rngSprint = $ws.Range("F2:F4")
$rngSprint.Validation.Delete()
$rngSprint.Validation.Add(3, 1, 1, '"Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10"')

# --- Selection cursor moved ---
$ws.Range("D7").Select()

Write-Host "done"
